$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "Handed back: in sync with en-US"

# --- Latest Handback DateTime refreshed for both languages ---
$zhcn.Range("K2").Value = "2016-10-19 11:11:45"
$dede.Range("K2").Value = "2016-10-19 11:12:02"

# --- Error Detail cleared now that handback versions are current ---
$zhcn.Range("P2").Value = ""
$dede.Range("P2").Value = ""

# --- Column width adjustments (widened to fit the longer status text) ---
$overview.Columns.Item(5).ColumnWidth = 29.14437166849777
$overview.Columns.Item(6).ColumnWidth = 29.14437166849777

$zhcn.Columns.Item(3).ColumnWidth = 29.14437166849777
$zhcn.Columns.Item(16).ColumnWidth = 12.913719813028965

$dede.Columns.Item(3).ColumnWidth = 29.14437166849777
$dede.Columns.Item(16).ColumnWidth = 12.913719813028965
